# Applies the cryptos.xlsx price/volume update described in the commit message:
# "Updated cryptos list on Tue Oct 10 02:33:52 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.667.85"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.583.92"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("D5").Value = "'206.95"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "'22.11"
$ws.Range("E8").Value = "  -4.79%  "
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "1.809.02"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").Value = "1.599.48"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("E14").Value = "  -4.13%  "
$ws.Range("E15").Value = "  -5.18%  "
$ws.Range("D16").Value = "'63.46"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "27.630.05"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'218.90"
$ws.Range("E18").Value = "  -3.98%  "
$ws.Range("D19").Value = "0.0₃0693"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").Value = "'7.30"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").Value = "'153.79"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "'6.84"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "'15.09"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("D33").Value = "1.361.36"
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D36").Value = "'0.968"
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("D37").Value = "'2.30"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("D40").Value = "'0.818"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "'0.969"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "'63.65"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.19"
$ws.Range("E45").Value = "  -3.99%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.73"
$ws.Range("E46").Value = "  -4.47%  "
$ws.Range("D47").Value = "1.719.99"
$ws.Range("D48").Value = "'88.05"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  +10.40%  "
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("E51").Value = "  -1.16%  "
